# Expense Req latest changes
$wb = $excel.ActiveWorkbook

$usersSheet = $wb.Worksheets.Item("Users")

# Delete the row containing "Amanda Donovan" (row 3), shifting "Leslie Ward" up to row 3.
$usersSheet.Rows.Item(3).Delete()

# Update the selection on the Users sheet.
$usersSheet.Range("D8").Select()

# Make "Users" the active/selected sheet (tab) in the workbook.
$usersSheet.Activate()
